$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the shared-string values in B2:B6 to drop the trailing noise suffix
$ws.Range("B2").Value = "S14075533.HN66-E4.1419414334"
$ws.Range("B3").Value = "S14075533.MB9-16-B1.1097600185"
$ws.Range("B4").Value = "S14075533.MB9-04-B6.1267956153"
$ws.Range("B5").Value = "S14075533.MB3-06-G25.1721880608"
$ws.Range("B6").Value = "S14075533.MB25-30-R21.1350325011"

# Move the active selection on Sheet1 from E6 to H6
$ws.Range("H6").Select()

# Scroll the workbook window back to the top-left (xWindow 28680 -> -120)
$excel.ActiveWindow.WindowState = -4143
